# Scheduled-runner price/profit refresh for the Pandaemonium Leve Profits
# workbook. Updates currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) on the affected leve rows across the per-job sheets, using the
# latest market-board pull. Columns A:G (item/leve metadata) are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1448.4615
$ws.Range("I112").Value = 300
$ws.Range("J112").Value = 1958.8889
$ws.Range("K112").Value = 900
$ws.Range("L112").Value = 5876.6667
$ws.Range("M112").Value = 208
$ws.Range("N112").Value = -8092.6667

$ws.Range("H137").Value = 2822.6897
$ws.Range("I137").Value = 2193.75
$ws.Range("J137").Value = 4220.3335
$ws.Range("K137").Value = 6581.25
$ws.Range("L137").Value = 12661.0005
$ws.Range("M137").Value = -4031.25
$ws.Range("N137").Value = -17761.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8136.886
$ws.Range("I61").Value = 4559.5625
$ws.Range("K61").Value = 4559.5625
$ws.Range("M61").Value = -4347.5625

$ws.Range("H74").Value = 4095
$ws.Range("I74").Value = 1918.4412
$ws.Range("J74").Value = 12317.556
$ws.Range("K74").Value = 1918.4412
$ws.Range("L74").Value = 12317.556
$ws.Range("M74").Value = -1044.4412
$ws.Range("N74").Value = -14065.556

$ws.Range("H77").Value = 4095
$ws.Range("I77").Value = 1918.4412
$ws.Range("J77").Value = 12317.556
$ws.Range("K77").Value = 9592.206
$ws.Range("L77").Value = 61587.78
$ws.Range("M77").Value = -5224.206
$ws.Range("N77").Value = -70323.78

$ws.Range("H132").Value = 2739.5557
$ws.Range("I132").Value = 1997.4117
$ws.Range("J132").Value = 4001.2
$ws.Range("K132").Value = 5992.2351
$ws.Range("L132").Value = 12003.6
$ws.Range("M132").Value = -3462.2351
$ws.Range("N132").Value = -17063.6

$ws.Range("H136").Value = 8136.886
$ws.Range("I136").Value = 4559.5625
$ws.Range("K136").Value = 13678.6875
$ws.Range("M136").Value = -11128.6875

$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 10026
$ws.Range("I32").Value = 10026
$ws.Range("K32").Value = 10026
$ws.Range("M32").Value = -9642

$ws.Range("H134").Value = 42744.56
$ws.Range("I134").Value = 3345.077
$ws.Range("K134").Value = 10035.231
$ws.Range("M134").Value = -7500.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4170.391
$ws.Range("I31").Value = 6172.2
$ws.Range("J31").Value = 2630.5386
$ws.Range("K31").Value = 6172.2
$ws.Range("L31").Value = 2630.5386
$ws.Range("M31").Value = -5877.2
$ws.Range("N31").Value = -3220.5386

$ws.Range("H34").Value = 4170.391
$ws.Range("I34").Value = 6172.2
$ws.Range("J34").Value = 2630.5386
$ws.Range("K34").Value = 6172.2
$ws.Range("L34").Value = 2630.5386
$ws.Range("M34").Value = -5970.2
$ws.Range("N34").Value = -3034.5386

$ws.Range("H50").Value = 13577.111
$ws.Range("J50").Value = 13577.111
$ws.Range("L50").Value = 13577.111
$ws.Range("N50").Value = -14827.111

$ws.Range("H58").Value = 6063307
$ws.Range("I58").Value = 9092837
$ws.Range("J58").Value = 4248
$ws.Range("K58").Value = 9092837
$ws.Range("L58").Value = 4248
$ws.Range("M58").Value = -9092634
$ws.Range("N58").Value = -4654

$ws.Range("H132").Value = 2747
$ws.Range("I132").Value = 2548.75
$ws.Range("J132").Value = 3086.8572
$ws.Range("K132").Value = 7646.25
$ws.Range("L132").Value = 9260.5716
$ws.Range("M132").Value = -5116.25
$ws.Range("N132").Value = -14320.5716

$ws.Range("H134").Value = 3006.4075
$ws.Range("J134").Value = 3023.7778
$ws.Range("L134").Value = 9071.3334
$ws.Range("N134").Value = -14141.3334

$ws.Range("H136").Value = 6063307
$ws.Range("I136").Value = 9092837
$ws.Range("J136").Value = 4248
$ws.Range("K136").Value = 27278511
$ws.Range("L136").Value = 12744
$ws.Range("M136").Value = -27275961
$ws.Range("N136").Value = -17844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 5000
$ws.Range("J93").Value = 5000
$ws.Range("L93").Value = 15000
$ws.Range("N93").Value = -18744

$ws.Range("H132").Value = 2580.5
$ws.Range("I132").Value = 5400
$ws.Range("J132").Value = 1875.625
$ws.Range("K132").Value = 48600
$ws.Range("L132").Value = 16880.625
$ws.Range("M132").Value = -46070
$ws.Range("N132").Value = -21940.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2633.3333
$ws.Range("I126").Value = 1920
$ws.Range("J126").Value = 3142.8572
$ws.Range("K126").Value = 5760
$ws.Range("L126").Value = 9428.5716
$ws.Range("M126").Value = -3290
$ws.Range("N126").Value = -14368.5716

$ws.Range("H132").Value = 5930.2583
$ws.Range("I132").Value = 2281.7144
$ws.Range("J132").Value = 13592.2
$ws.Range("K132").Value = 6845.1432
$ws.Range("L132").Value = 40776.60000000001
$ws.Range("M132").Value = -4315.1432
$ws.Range("N132").Value = -45836.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5288.8823
$ws.Range("I132").Value = 4916.3076
$ws.Range("J132").Value = 6499.75
$ws.Range("K132").Value = 14748.9228
$ws.Range("L132").Value = 19499.25
$ws.Range("M132").Value = -12218.9228
$ws.Range("N132").Value = -24559.25

$ws.Range("H136").Value = 7443.4346
$ws.Range("I136").Value = 4743.9
$ws.Range("J136").Value = 9520
$ws.Range("K136").Value = 14231.7
$ws.Range("L136").Value = 28560
$ws.Range("M136").Value = -11681.7
$ws.Range("N136").Value = -33660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 45662.5
$ws.Range("J86").Value = 45662.5
$ws.Range("L86").Value = 45662.5
$ws.Range("N86").Value = -47908.5

$ws.Range("H89").Value = 45662.5
$ws.Range("J89").Value = 45662.5
$ws.Range("L89").Value = 228312.5
$ws.Range("N89").Value = -239544.5

$ws.Range("H132").Value = 1711.6086
$ws.Range("I132").Value = 1045.0312
$ws.Range("J132").Value = 3235.2144
$ws.Range("K132").Value = 3135.0936
$ws.Range("L132").Value = 9705.643199999999
$ws.Range("M132").Value = -605.0935999999997
$ws.Range("N132").Value = -14765.6432

$ws.Range("H136").Value = 5935.4614
$ws.Range("I136").Value = 2580.9473
$ws.Range("J136").Value = 9122.25
$ws.Range("K136").Value = 7742.841899999999
$ws.Range("L136").Value = 27366.75
$ws.Range("M136").Value = -5192.841899999999
$ws.Range("N136").Value = -32466.75
